{"js": "// Apply the four textual edits described by the diff:\n//  1. Insert \"AL \" right before the \"{{NUMERO_VENDEDOR}}\" placeholder.\n//  2. Drop the curly quotes around RELEVANTES (\"RELEVANTES\" -> RELEVANTES).\n//  3. Replace the literal word \"COMPRADOR\" (in \"PROMITENTE COMPRADOR\")\n//     with the \"{{SEXO_4}}\" placeholder, keeping the separating space.\n//  4. Replace the literal word \"VENDEDOR\" (in \"PROMITENTE VENDEDOR\")\n//     with the \"{{SEXO_2}}\" placeholder, keeping the separating space.\n\nconst body = context.document.body;\n\n// ---- 1. \"...NOTIFICACIONES {{NUMERO_VENDEDOR}}...\" -> \"...NOTIFICACIONES AL {{NUMERO_VENDEDOR}}...\"\nconst numeroVendedorHits = body.search(\"{{NUMERO_VENDEDOR}}\", { matchCase: true });\nnumeroVendedorHits.load(\"items\");\nawait context.sync();\nif (numeroVendedorHits.items.length > 0) {\n  numeroVendedorHits.items[0].insertText(\"AL \", Word.InsertLocation.before);\n}\nawait context.sync();\n\n// ---- 2. Remove the curly quotes wrapped around RELEVANTES.\nconst relevantesHits = body.search(\"\\u201cRELEVANTES\\u201d\", { matchCase: true });\nrelevantesHits.load(\"items\");\nawait context.sync();\nif (relevantesHits.items.length > 0) {\n  relevantesHits.items[0].insertText(\"RELEVANTES\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// ---- 3. \"PROMITENTE COMPRADOR\" -> \"PROMITENTE {{SEXO_4}}\"\nconst compradorAnchor = body.search(\"PROMITENTE COMPRADOR\", { matchCase: true });\ncompradorAnchor.load(\"items\");\nawait context.sync();\nif (compradorAnchor.items.length > 0) {\n  const compradorSub = compradorAnchor.items[0].search(\" COMPRADOR\", { matchCase: true });\n  compradorSub.load(\"items\");\n  await context.sync();\n  if (compradorSub.items.length > 0) {\n    compradorSub.items[0].insertText(\" {{SEXO_4}}\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// ---- 4. \"PROMITENTE VENDEDOR\" (the signature-block one, next to \"PROMITENTE COMPRADOR\"\n//         a.k.a. now \"PROMITENTE {{SEXO_4}}\") -> \"PROMITENTE {{SEXO_2}}\"\nconst vendedorAnchors = body.search(\"PROMITENTE VENDEDOR\", { matchCase: true });\nvendedorAnchors.load(\"items\");\nawait context.sync();\n\nlet vendedorTarget = null;\nfor (const item of vendedorAnchors.items) {\n  const para = item.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n  // The signature line is the only paragraph that also carries (or used to\n  // carry, before step 3 above ran) \"COMPRADOR\" right next to \"PROMITENTE\".\n  if (para.text.indexOf(\"COMPRADOR\") !== -1 || para.text.indexOf(\"{{SEXO_4}}\") !== -1) {\n    vendedorTarget = item;\n    break;\n  }\n}\n\nif (vendedorTarget) {\n  const vendedorSub = vendedorTarget.search(\" VENDEDOR\", { matchCase: true });\n  vendedorSub.load(\"items\");\n  await context.sync();\n  if (vendedorSub.items.length > 0) {\n    vendedorSub.items[0].insertText(\" {{SEXO_2}}\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the four textual edits described by the diff:\n#  1. Insert \"AL \" right before the \"{{NUMERO_VENDEDOR}}\" placeholder.\n#  2. Drop the curly quotes around RELEVANTES (\"RELEVANTES\" -> RELEVANTES).\n#  3. Replace the literal word \"COMPRADOR\" (in \"PROMITENTE COMPRADOR\")\n#     with the \"{{SEXO_4}}\" placeholder, keeping the separating space.\n#  4. Replace the literal word \"VENDEDOR\" (in \"PROMITENTE VENDEDOR\")\n#     with the \"{{SEXO_2}}\" placeholder, keeping the separating space.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne   = 2\n$openCurly  = [char]0x201C\n$closeCurly = [char]0x201D\n\n# ---- 1. \"...NOTIFICACIONES {{NUMERO_VENDEDOR}}...\" -> \"...NOTIFICACIONES AL {{NUMERO_VENDEDOR}}...\"\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Replacement.ClearFormatting()\n$r1.Find.Text = \"{{NUMERO_VENDEDOR}}\"\n$r1.Find.Replacement.Text = \"AL {{NUMERO_VENDEDOR}}\"\n$r1.Find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceOne)\n\n# ---- 2. Remove the curly quotes wrapped around RELEVANTES.\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Replacement.ClearFormatting()\n$r2.Find.Text = $openCurly + \"RELEVANTES\" + $closeCurly\n$r2.Find.Replacement.Text = \"RELEVANTES\"\n$r2.Find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceOne)\n\n# ---- 3 & 4: the signature line \"PROMITENTE COMPRADOR ... PROMITENTE VENDEDOR\".\n# Locate that specific paragraph first so the COMPRADOR/VENDEDOR replacements\n# cannot accidentally hit the unrelated \"...QUE EL PROMITENTE VENDEDOR LE\n# ENTREGUE...\" sentence elsewhere in the document.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*PROMITENTE COMPRADOR*PROMITENTE VENDEDOR*\") {\n    $rComprador = $p.Range\n    $rComprador.Find.ClearFormatting()\n    $rComprador.Find.Replacement.ClearFormatting()\n    $rComprador.Find.Text = \" COMPRADOR\"\n    $rComprador.Find.Replacement.Text = \" {{SEXO_4}}\"\n    $rComprador.Find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceOne)\n\n    $rVendedor = $p.Range\n    $rVendedor.Find.ClearFormatting()\n    $rVendedor.Find.Replacement.ClearFormatting()\n    $rVendedor.Find.Text = \" VENDEDOR\"\n    $rVendedor.Find.Replacement.Text = \" {{SEXO_2}}\"\n    $rVendedor.Find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceOne)\n\n    break\n  }\n}\n"}
